$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name to reflect new "through" date
$ws.Name = "Through 2022-09-19"

# Update the "September (through 09-17)" label to "September (through 09-19)"
$ws.Range("A10").Value = "September (through 09-19)"

# Update September row (row 10) values
$ws.Range("B10").Value = 21
$ws.Range("C10").Value = 34
$ws.Range("D10").Value = 42
$ws.Range("E10").Value = 36
$ws.Range("F10").Value = 44
$ws.Range("G10").Value = 74
$ws.Range("H10").Value = 109
$ws.Range("I10").Value = 92

# Update Total row (row 11) values
$ws.Range("B11").Value = 215
$ws.Range("C11").Value = 415
$ws.Range("D11").Value = 593
$ws.Range("E11").Value = 526
$ws.Range("F11").Value = 393
$ws.Range("G11").Value = 858
$ws.Range("H11").Value = 1179
$ws.Range("I11").Value = 1227
